$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.06809999999999
$ws.Range("B9").Value = 8.574100000000005
$ws.Range("B18").Value = 4.833800000000005
$ws.Range("B20").Value = 5.631899999999998
$ws.Range("E21").Value = 13.14589999999999
